$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 479.83334
$ws.Range("I2").Value = 261.77777
$ws.Range("K2").Value = 261.77777
$ws.Range("M2").Value = -148.77777
$ws.Range("H43").Value = 7167.1665
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").ClearContents()
$ws.Range("N48").Value = 0
$ws.Range("H51").Value = 15000
$ws.Range("I51").Value = 15000
$ws.Range("K51").Value = 15000
$ws.Range("M51").Value = -14516
$ws.Range("H55").Value = 124.583336
$ws.Range("I55").Value = 71.57143000000001
$ws.Range("J55").Value = 198.8
$ws.Range("K55").Value = 71.57143000000001
$ws.Range("L55").Value = 198.8
$ws.Range("M55").Value = 142.42857
$ws.Range("N55").Value = -626.8
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").ClearContents()
$ws.Range("N56").Value = 0
$ws.Range("H80").Value = 3567
$ws.Range("I80").Value = 3233.3333
$ws.Range("J80").Value = 3733.8333
$ws.Range("K80").Value = 9699.999899999999
$ws.Range("L80").Value = 11201.4999
$ws.Range("M80").Value = -8701.999899999999
$ws.Range("N80").Value = -13197.4999
$ws.Range("H83").Value = 3567
$ws.Range("I83").Value = 3233.3333
$ws.Range("J83").Value = 3733.8333
$ws.Range("K83").Value = 29099.9997
$ws.Range("L83").Value = 33604.4997
$ws.Range("M83").Value = -24107.9997
$ws.Range("N83").Value = -43588.4997
$ws.Range("H135").Value = 1650.8572
$ws.Range("I135").Value = 1742.6666
$ws.Range("K135").Value = 15683.9994
$ws.Range("M135").Value = -13148.9994
$ws.Range("H137").Value = 3561.875
$ws.Range("I137").Value = 1996
$ws.Range("J137").Value = 3785.5715
$ws.Range("K137").Value = 5988
$ws.Range("L137").Value = 11356.7145
$ws.Range("M137").Value = -3438
$ws.Range("N137").Value = -16456.7145
$ws.Range("H141").Value = 500
$ws.Range("I141").Value = 500
$ws.Range("K141").Value = 1500
$ws.Range("M141").Value = 3680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1805.4546
$ws.Range("I45").Value = 1762.5555
$ws.Range("J45").Value = 1998.5
$ws.Range("K45").Value = 1762.5555
$ws.Range("L45").Value = 1998.5
$ws.Range("M45").Value = -1385.5555
$ws.Range("N45").Value = -2752.5
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0
$ws.Range("H110").Value = 1168.3334
$ws.Range("I110").Value = 1074.2858
$ws.Range("J110").Value = 1497.5
$ws.Range("K110").Value = 1074.2858
$ws.Range("L110").Value = 1497.5
$ws.Range("M110").Value = 970.7141999999999
$ws.Range("N110").Value = -5587.5
$ws.Range("H132").Value = 2505.5
$ws.Range("I132").Value = 1793.0667
$ws.Range("J132").Value = 4642.8
$ws.Range("K132").Value = 5379.2001
$ws.Range("L132").Value = 13928.4
$ws.Range("M132").Value = -2849.2001
$ws.Range("N132").Value = -18988.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 14800
$ws.Range("J76").Value = 14800
$ws.Range("L76").Value = 14800
$ws.Range("N76").Value = -15430
$ws.Range("H79").Value = 14800
$ws.Range("J79").Value = 14800
$ws.Range("L79").Value = 14800
$ws.Range("N79").Value = -16984
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").ClearContents()
$ws.Range("N114").Value = 0
$ws.Range("H134").Value = 1598.0526
$ws.Range("I134").Value = 1409.0555
$ws.Range("K134").Value = 4227.166499999999
$ws.Range("M134").Value = -1692.166499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2664.3572
$ws.Range("I31").Value = 2788
$ws.Range("K31").Value = 2788
$ws.Range("M31").Value = -2493
$ws.Range("H34").Value = 2664.3572
$ws.Range("I34").Value = 2788
$ws.Range("K34").Value = 2788
$ws.Range("M34").Value = -2586
$ws.Range("H58").Value = 12002.2
$ws.Range("J58").Value = 12502.75
$ws.Range("L58").Value = 12502.75
$ws.Range("N58").Value = -12908.75
$ws.Range("H60").Value = 11353.333
$ws.Range("I60").Value = 2030
$ws.Range("J60").Value = 30000
$ws.Range("K60").Value = 2030
$ws.Range("L60").Value = 30000
$ws.Range("M60").Value = -1519
$ws.Range("N60").Value = -31022
$ws.Range("H136").Value = 12002.2
$ws.Range("J136").Value = 12502.75
$ws.Range("L136").Value = 37508.25
$ws.Range("N136").Value = -42608.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 167253.22
$ws.Range("J12").Value = 1023.3
$ws.Range("L12").Value = 3069.9
$ws.Range("N12").Value = -3415.9
$ws.Range("H131").Value = 1729.766
$ws.Range("I131").Value = 1649.5
$ws.Range("J131").Value = 1733.3334
$ws.Range("K131").Value = 4948.5
$ws.Range("L131").Value = 5200.0002
$ws.Range("M131").Value = 91.5
$ws.Range("N131").Value = -15280.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 60000
$ws.Range("I87").Value = 60000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 60000
$ws.Range("L87").ClearContents()
$ws.Range("M87").Value = -58752
$ws.Range("N87").Value = 0
$ws.Range("H90").Value = 60000
$ws.Range("I90").Value = 60000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 180000
$ws.Range("L90").ClearContents()
$ws.Range("M90").Value = -173760
$ws.Range("N90").Value = 0
$ws.Range("H102").Value = 638.75
$ws.Range("I102").Value = 638.75
$ws.Range("K102").Value = 638.75
$ws.Range("M102").Value = 983.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 15000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15386
$ws.Range("H122").Value = 2277.3333
$ws.Range("J122").Value = 2132.6667
$ws.Range("L122").Value = 6398.000100000001
$ws.Range("N122").Value = -11298.0001
$ws.Range("H136").Value = 4122.846
$ws.Range("I136").Value = 3790.6365
$ws.Range("K136").Value = 11371.9095
$ws.Range("M136").Value = -8821.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 34750
$ws.Range("J82").Value = 34750
$ws.Range("L82").Value = 34750
$ws.Range("N82").Value = -35516
$ws.Range("H85").Value = 34750
$ws.Range("J85").Value = 34750
$ws.Range("L85").Value = 34750
